# Insert two new weekly records for "Betarraga" (Hortaliza) right before the
# existing row 242 block, pushing the remaining rows (old 242-273) down to
# 244-275. The two new rows replicate the row immediately below them (so all
# the invariant columns - Mercado, Region, Categoria, etc. - come along for
# free) and then the columns that actually differ for the new records are
# overwritten explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert two blank rows at 242 (each Insert() pushes everything
# at/after the target row down by one).
$ws.Rows.Item(242).Insert()
$ws.Rows.Item(242).Insert()

# Seed the two new rows with the contents of the rows that are now directly
# beneath them (the original row 242 data, now living at rows 244/245) so
# every column starts out populated, then we only need to patch the columns
# that actually change.
$lastCol = 18
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item(242, $c).Value2 = $ws.Cells.Item(244, $c).Value2
    $ws.Cells.Item(243, $c).Value2 = $ws.Cells.Item(245, $c).Value2
}

# New row 242: Primera, 22-Oct-2021 (serial 44491)
$ws.Cells.Item(242, 4).Value2 = 44491
$ws.Cells.Item(242, 9).Value2 = "Primera"
$ws.Cells.Item(242, 10).Value2 = 5200
$ws.Cells.Item(242, 11).Value2 = 90
$ws.Cells.Item(242, 12).Value2 = 100
$ws.Cells.Item(242, 13).Value2 = 95
$ws.Cells.Item(242, 16).Value2 = 95

# New row 243: Segunda, 22-Oct-2021 (serial 44491)
$ws.Cells.Item(243, 4).Value2 = 44491
$ws.Cells.Item(243, 9).Value2 = "Segunda"
$ws.Cells.Item(243, 10).Value2 = 2500
$ws.Cells.Item(243, 11).Value2 = 70
$ws.Cells.Item(243, 12).Value2 = 80
$ws.Cells.Item(243, 13).Value2 = 75
$ws.Cells.Item(243, 16).Value2 = 75
